$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "MEC-3A-Trat. Termicos" entry from column B (segunda) to column D (quarta)
# for rows 11 (13:50) and 12 (14:40), matching the new schedule.

$ws.Range("B11").Value = "-"
$ws.Range("D11").Value = "[-, 'MEC-3A-Trat. Termicos', -, -]"

$ws.Range("B12").Value = "-"
$ws.Range("D12").Value = "[-, 'MEC-3A-Trat. Termicos', -, -]"
